$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Incomplete" labels from column E (they are dropped entirely in the edit)
$ws.Range("E7").Value = $null
$ws.Range("E11").Value = $null

# Insert a new row before the current row 10 ("Bayesian Optimisation for XGBOOST ...")
# so that a "Reticulate" entry moves up to sit right after the "PSO optimisation" row.
$ws.Rows.Item(10).Insert()

# The newly inserted row 10 becomes "Reticulate" (previously the last row of the list)
$ws.Range("B10").Value = "Reticulate"

# Clear out the now-duplicate "Reticulate" entry that got pushed down to row 16,
# and replace it with the new "Loan portfolio optimisation" entry / series info.
$ws.Range("B16").Value = "Loan portfolio optimisation"
$ws.Range("C16").Value = "Credit risk series"

# Update the used range / selection to reflect the new last edited cell
$ws.Range("D16").Select()
